$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 20, pushing the existing rows 20-37 down to 21-38.
$ws.Rows(20).Insert()

# Populate the newly inserted row 20 with the new weekly record.
$ws.Range("A20").Value = 7
$ws.Range("B20").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C20").Value = "Ñuble"
$ws.Range("D20").Value = 44601
$ws.Range("E20").Value = 16
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100103
$ws.Range("H20").Value = "Frutos de hueso (carozo)"
$ws.Range("I20").Value = 100103002
$ws.Range("J20").Value = "Ciruela"
$ws.Range("K20").Value = "Black Amber"
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 60
$ws.Range("N20").Value = 9000
$ws.Range("O20").Value = 10000
$ws.Range("P20").Value = 9500
$ws.Range("Q20").Value = "$/bandeja 18 kilos granel"
$ws.Range("R20").Value = "Provincia de Curicó"
$ws.Range("S20").Value = 528
$ws.Range("T20").Value = 18
